$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume updates (GitHub Actions sync, 2023-01-08 03:56 UTC).
# D (Price) and E (Volume(1h)) are plain-text cells holding numeric-looking
# strings; a leading apostrophe stops Excel from coercing them to Number cells,
# and resetting the Style back to "Normal" afterwards clears the quote-prefix
# formatting it implicitly applies, so the cell keeps its original (default) style.

$ws.Range("D2").Value = "'260.84"
$ws.Range("E2").Value = "'-0.38%"
$ws.Range("D2:E2").Style = "Normal"

$ws.Range("D3").Value = "'27.08"
$ws.Range("E3").Value = "'-1.09%"
$ws.Range("D3:E3").Style = "Normal"

$ws.Range("D4").Value = "'4.698"
$ws.Range("E4").Value = "'-1.15%"
$ws.Range("D4:E4").Style = "Normal"

$ws.Range("E5").Value = "'2.47%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'6.731"
$ws.Range("E6").Value = "'0.18%"
$ws.Range("D6:E6").Style = "Normal"

$ws.Range("D7").Value = "'0.8511"
$ws.Range("E7").Value = "'-1.74%"
$ws.Range("D7:E7").Style = "Normal"

$ws.Range("D8").Value = "'0.9095"
$ws.Range("E8").Value = "'-1.70%"
$ws.Range("D8:E8").Style = "Normal"

$ws.Range("D9").Value = "'0.1400"
$ws.Range("E9").Value = "'-0.85%"
$ws.Range("D9:E9").Style = "Normal"

$ws.Range("D10").Value = "'0.04878"
$ws.Range("E10").Value = "'-2.31%"
$ws.Range("D10:E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07083"
$ws.Range("E11").Value = "'-1.05%"
$ws.Range("D11:E11").Style = "Normal"

$ws.Range("D12").Value = "'0.03128"
$ws.Range("E12").Value = "'2.88%"
$ws.Range("D12:E12").Style = "Normal"

$ws.Range("D13").Value = "'0.09057"
$ws.Range("E13").Value = "'-0.65%"
$ws.Range("D13:E13").Style = "Normal"

$ws.Range("E14").Value = "'-0.39%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.0006145"
$ws.Range("E15").Value = "'0.95%"
$ws.Range("D15:E15").Style = "Normal"

$ws.Range("D16").Value = "'0.006086"
$ws.Range("E16").Value = "'-1.76%"
$ws.Range("D16:E16").Style = "Normal"

$ws.Range("D17").Value = "'3.453"
$ws.Range("E17").Value = "'0.14%"
$ws.Range("D17:E17").Style = "Normal"

$ws.Range("D18").Value = "'3.171"
$ws.Range("E18").Value = "'0.10%"
$ws.Range("D18:E18").Style = "Normal"

$ws.Range("E21").Value = "'0.84%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'4.104"
$ws.Range("E22").Value = "'0.28%"
$ws.Range("D22:E22").Style = "Normal"

$ws.Range("D23").Value = "'0.04252"
$ws.Range("E23").Value = "'-0.09%"
$ws.Range("D23:E23").Style = "Normal"

$ws.Range("D24").Value = "'0.001219"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").Value = "'0.004081"
$ws.Range("E25").Value = "'4.31%"
$ws.Range("D25:E25").Style = "Normal"

$ws.Range("D26").Value = "'0.0001200"
$ws.Range("E26").Value = "'0.03%"
$ws.Range("D26:E26").Style = "Normal"

$ws.Range("D27").Value = "'0.0001639"
$ws.Range("E27").Value = "'4.32%"
$ws.Range("D27:E27").Style = "Normal"

$ws.Range("D40").Value = "'0.03941"
$ws.Range("E40").Value = "'1.49%"
$ws.Range("D40:E40").Style = "Normal"

$ws.Range("D41").Value = "'0.1110"
$ws.Range("E41").Value = "'-0.46%"
$ws.Range("D41:E41").Style = "Normal"

$ws.Range("D42").Value = "'0.004126"
$ws.Range("E42").Value = "'-0.16%"
$ws.Range("D42:E42").Style = "Normal"

$ws.Range("D43").Value = "'0.002151"
$ws.Range("E43").Value = "'-2.58%"
$ws.Range("D43:E43").Style = "Normal"

$ws.Range("E44").Value = "'-7.19%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.00005124"
$ws.Range("E45").Value = "'-4.48%"
$ws.Range("D45:E45").Style = "Normal"

$ws.Range("E46").Value = "'-0.03%"
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.03400"
$ws.Range("E47").Value = "'-37.69%"
$ws.Range("D47:E47").Style = "Normal"

$ws.Range("D48").Value = "'0.06889"
$ws.Range("E48").Value = "'-49.08%"
$ws.Range("D48:E48").Style = "Normal"

$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'-0.03%"
$ws.Range("D49:E49").Style = "Normal"

$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'-0.03%"
$ws.Range("D50:E50").Style = "Normal"
